# Case_5_140 (380 kV) res_line/pl_mw.xlsx results refresh.
# Updates the computed columns (B-F, I, L-N) for rows 2-25 (A = 0..23)
# to the new simulation output values; G, H, J, K, O stay 0 (unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 0)
$ws.Cells.Item(2, 2).Value = 1.775537327644599
$ws.Cells.Item(2, 3).Value = 0.2726853044468669
$ws.Cells.Item(2, 4).Value = 0.08119606441675487
$ws.Cells.Item(2, 5).Value = 0.05107260235544597
$ws.Cells.Item(2, 6).Value = 1.904377458661031
$ws.Cells.Item(2, 9).Value = 1.389428456148138
$ws.Cells.Item(2, 12).Value = 0.2189841073584375
$ws.Cells.Item(2, 13).Value = 0.3570734750569073
$ws.Cells.Item(2, 14).Value = 1.869836136567365

# Row 3 (A3 = 1)
$ws.Cells.Item(3, 2).Value = 1.666964319022725
$ws.Cells.Item(3, 3).Value = 0.237957907749319
$ws.Cells.Item(3, 4).Value = 0.0820589068039741
$ws.Cells.Item(3, 5).Value = 0.05048635519244993
$ws.Cells.Item(3, 6).Value = 1.874001328707891
$ws.Cells.Item(3, 9).Value = 1.384372155566005
$ws.Cells.Item(3, 12).Value = 0.2161868950077377
$ws.Cells.Item(3, 13).Value = 0.3408409847153564
$ws.Cells.Item(3, 14).Value = 1.889697917550571

# Row 4 (A4 = 2)
$ws.Cells.Item(4, 2).Value = 1.601271572767757
$ws.Cells.Item(4, 3).Value = 0.2166670114999931
$ws.Cells.Item(4, 4).Value = 0.08261722630306423
$ws.Cells.Item(4, 5).Value = 0.05011974870084313
$ws.Cells.Item(4, 6).Value = 1.856591702773073
$ws.Cells.Item(4, 9).Value = 1.38206801549984
$ws.Cells.Item(4, 12).Value = 0.2145839486543437
$ws.Cells.Item(4, 13).Value = 0.3310786735456759
$ws.Cells.Item(4, 14).Value = 1.902560060508286

# Row 5 (A5 = 3)
$ws.Cells.Item(5, 2).Value = 1.574744687005307
$ws.Cells.Item(5, 3).Value = 0.2079981795276922
$ws.Cells.Item(5, 4).Value = 0.0828519066576181
$ws.Cells.Item(5, 5).Value = 0.04996866713594539
$ws.Cells.Item(5, 6).Value = 1.849807918807059
$ws.Cells.Item(5, 9).Value = 1.381329678822837
$ws.Cells.Item(5, 12).Value = 0.2139595162056267
$ws.Cells.Item(5, 13).Value = 0.3271517463758897
$ws.Cells.Item(5, 14).Value = 1.907968852603943

# Row 6 (A6 = 4)
$ws.Cells.Item(6, 2).Value = 1.570354596775928
$ws.Cells.Item(6, 3).Value = 0.2065591540626883
$ws.Cells.Item(6, 4).Value = 0.08289130724132576
$ws.Cells.Item(6, 5).Value = 0.04994347787413567
$ws.Cells.Item(6, 6).Value = 1.848700206632486
$ws.Cells.Item(6, 9).Value = 1.381219176311511
$ws.Cells.Item(6, 12).Value = 0.2138575673135321
$ws.Cells.Item(6, 13).Value = 0.3265027774379092
$ws.Cells.Item(6, 14).Value = 1.908877076931738

# Row 7 (A7 = 5)
$ws.Cells.Item(7, 2).Value = 1.600912837343628
$ws.Cells.Item(7, 3).Value = 0.2165500715386486
$ws.Cells.Item(7, 4).Value = 0.08262036230929404
$ws.Cells.Item(7, 5).Value = 0.05011771800912701
$ws.Cells.Item(7, 6).Value = 1.856498957802486
$ws.Cells.Item(7, 9).Value = 1.382057246618885
$ws.Cells.Item(7, 12).Value = 0.2145754108366305
$ws.Cells.Item(7, 13).Value = 0.3310255061156369
$ws.Cells.Item(7, 14).Value = 1.90263232818549

# Row 8 (A8 = 6)
$ws.Cells.Item(8, 2).Value = 1.737899156216145
$ws.Cells.Item(8, 3).Value = 0.2607042924405221
$ws.Cells.Item(8, 4).Value = 0.0814876426940554
$ws.Cells.Item(8, 5).Value = 0.05087183522592653
$ws.Cells.Item(8, 6).Value = 1.893645253106669
$ws.Cells.Item(8, 9).Value = 1.387518506095006
$ws.Cells.Item(8, 12).Value = 0.2179958366475958
$ws.Cells.Item(8, 13).Value = 0.3514339822397119
$ws.Cells.Item(8, 14).Value = 1.876545853849045

# Row 9 (A9 = 7)
$ws.Cells.Item(9, 2).Value = 2.014292829440024
$ws.Cells.Item(9, 3).Value = 0.3475762332720365
$ws.Cells.Item(9, 4).Value = 0.07949335182825301
$ws.Cells.Item(9, 5).Value = 0.05229861833653082
$ws.Cells.Item(9, 6).Value = 1.976408685374665
$ws.Cells.Item(9, 9).Value = 1.404611551402894
$ws.Cells.Item(9, 12).Value = 0.2256139871158211
$ws.Cells.Item(9, 13).Value = 0.3930861333790787
$ws.Cells.Item(9, 14).Value = 1.830696734479361

# Row 10 (A10 = 8)
$ws.Cells.Item(10, 2).Value = 2.222193218714949
$ws.Cells.Item(10, 3).Value = 0.4116268267911209
$ws.Cells.Item(10, 4).Value = 0.0781672914163849
$ws.Cells.Item(10, 5).Value = 0.05331621465525238
$ws.Cells.Item(10, 6).Value = 2.043367976329051
$ws.Cells.Item(10, 9).Value = 1.421109830478997
$ws.Cells.Item(10, 12).Value = 0.2317696862633625
$ws.Cells.Item(10, 13).Value = 0.4246980877480908
$ws.Cells.Item(10, 14).Value = 1.800265418809268

# Row 11 (A11 = 9)
$ws.Cells.Item(11, 2).Value = 2.317845775415549
$ws.Cells.Item(11, 3).Value = 0.4408257229050605
$ws.Cells.Item(11, 4).Value = 0.07759441983867532
$ws.Cells.Item(11, 5).Value = 0.05377272276942691
$ws.Cells.Item(11, 6).Value = 2.075189197488555
$ws.Cells.Item(11, 9).Value = 1.429481678669987
$ws.Cells.Item(11, 12).Value = 0.2346922070794193
$ws.Cells.Item(11, 13).Value = 0.4393022136490856
$ws.Cells.Item(11, 14).Value = 1.787132196650511

# Row 12 (A12 = 10)
$ws.Cells.Item(12, 2).Value = 2.354223332982485
$ws.Cells.Item(12, 3).Value = 0.4518923121091802
$ws.Cells.Item(12, 4).Value = 0.07738187074691361
$ws.Cells.Item(12, 5).Value = 0.05394468751331338
$ws.Cells.Item(12, 6).Value = 2.087436555009475
$ws.Cells.Item(12, 9).Value = 1.432777340659172
$ws.Cells.Item(12, 12).Value = 0.2358165217105039
$ws.Cells.Item(12, 13).Value = 0.4448648105816559
$ws.Cells.Item(12, 14).Value = 1.782261500054588

# Row 13 (A13 = 11)
$ws.Cells.Item(13, 2).Value = 2.346381819429325
$ws.Cells.Item(13, 3).Value = 0.4495084860094494
$ws.Cells.Item(13, 4).Value = 0.07742745166445353
$ws.Cells.Item(13, 5).Value = 0.05390769182810651
$ws.Cells.Item(13, 6).Value = 2.084790064028482
$ws.Cells.Item(13, 9).Value = 1.432061969808913
$ws.Cells.Item(13, 12).Value = 0.2355735961837553
$ws.Cells.Item(13, 13).Value = 0.4436653651689895
$ws.Cells.Item(13, 14).Value = 1.783305922603809

# Row 14 (A14 = 12)
$ws.Cells.Item(14, 2).Value = 2.320835446167052
$ws.Cells.Item(14, 3).Value = 0.4417359814489146
$ws.Cells.Item(14, 4).Value = 0.07757684530443321
$ws.Cells.Item(14, 5).Value = 0.05378688845953228
$ws.Cells.Item(14, 6).Value = 2.076192830348106
$ws.Cells.Item(14, 9).Value = 1.429750296476712
$ws.Cells.Item(14, 12).Value = 0.2347843517767956
$ws.Cells.Item(14, 13).Value = 0.4397592029118442
$ws.Cells.Item(14, 14).Value = 1.786729423242072

# Row 15 (A15 = 13)
$ws.Cells.Item(15, 2).Value = 2.305207903739813
$ws.Cells.Item(15, 3).Value = 0.4369763689801403
$ws.Cells.Item(15, 4).Value = 0.07766892484906585
$ws.Cells.Item(15, 5).Value = 0.0537127755191662
$ws.Cells.Item(15, 6).Value = 2.070952526134676
$ws.Cells.Item(15, 9).Value = 1.428350687586942
$ws.Cells.Item(15, 12).Value = 0.234303212519734
$ws.Cells.Item(15, 13).Value = 0.4373707812380943
$ws.Cells.Item(15, 14).Value = 1.788839787395887

# Row 16 (A16 = 14)
$ws.Cells.Item(16, 2).Value = 2.215963872701252
$ws.Cells.Item(16, 3).Value = 0.4097199159260754
$ws.Cells.Item(16, 4).Value = 0.07820534247217203
$ws.Cells.Item(16, 5).Value = 0.05328625328719738
$ws.Cells.Item(16, 6).Value = 2.041315912135531
$ws.Cells.Item(16, 9).Value = 1.420580218374354
$ws.Cells.Item(16, 12).Value = 0.231581156119347
$ws.Cells.Item(16, 13).Value = 0.4237481889710253
$ws.Cells.Item(16, 14).Value = 1.80113803405483

# Row 17 (A17 = 15)
$ws.Cells.Item(17, 2).Value = 2.161492292326102
$ws.Cells.Item(17, 3).Value = 0.393015323111058
$ws.Cells.Item(17, 4).Value = 0.07854220674381196
$ws.Cells.Item(17, 5).Value = 0.05302296714042853
$ws.Cells.Item(17, 6).Value = 2.02348457130023
$ws.Cells.Item(17, 9).Value = 1.416035831370309
$ws.Cells.Item(17, 12).Value = 0.2299426037240266
$ws.Cells.Item(17, 13).Value = 0.4154485778649928
$ws.Cells.Item(17, 14).Value = 1.808864812694559

# Row 18 (A18 = 16)
$ws.Cells.Item(18, 2).Value = 2.130263046645496
$ws.Cells.Item(18, 3).Value = 0.3834130732057588
$ws.Cells.Item(18, 4).Value = 0.07873881860710874
$ws.Cells.Item(18, 5).Value = 0.05287092842616126
$ws.Cells.Item(18, 6).Value = 2.013356470684556
$ws.Cells.Item(18, 9).Value = 1.413503549759852
$ws.Cells.Item(18, 12).Value = 0.2290116572579137
$ws.Cells.Item(18, 13).Value = 0.4106959055676853
$ws.Cells.Item(18, 14).Value = 1.813375865228167

# Row 19 (A19 = 17)
$ws.Cells.Item(19, 2).Value = 2.119706739036985
$ws.Cells.Item(19, 3).Value = 0.3801628899207117
$ws.Cells.Item(19, 4).Value = 0.07880587802525341
$ws.Cells.Item(19, 5).Value = 0.05281934660947307
$ws.Cells.Item(19, 6).Value = 2.009949207314136
$ws.Cells.Item(19, 9).Value = 1.412660141456826
$ws.Cells.Item(19, 12).Value = 0.2286984297267622
$ws.Cells.Item(19, 13).Value = 0.4090903418561993
$ws.Cells.Item(19, 14).Value = 1.814914691756016

# Row 20 (A20 = 18)
$ws.Cells.Item(20, 2).Value = 2.167280386493644
$ws.Cells.Item(20, 3).Value = 0.3947929513054191
$ws.Cells.Item(20, 4).Value = 0.07850605119794452
$ws.Cells.Item(20, 5).Value = 0.05305105674657717
$ws.Cells.Item(20, 6).Value = 2.025369487194155
$ws.Cells.Item(20, 9).Value = 1.416511145401927
$ws.Cells.Item(20, 12).Value = 0.2301158393733118
$ws.Cells.Item(20, 13).Value = 0.4163299064315851
$ws.Cells.Item(20, 14).Value = 1.808035365686578

# Row 21 (A21 = 19)
$ws.Cells.Item(21, 2).Value = 2.328334792651333
$ws.Cells.Item(21, 3).Value = 0.4440186883259116
$ws.Cells.Item(21, 4).Value = 0.07753284560891949
$ws.Cells.Item(21, 5).Value = 0.0538223957729187
$ws.Cells.Item(21, 6).Value = 2.078712677140857
$ws.Cells.Item(21, 9).Value = 1.430425880590221
$ws.Cells.Item(21, 12).Value = 0.2350156934917607
$ws.Cells.Item(21, 13).Value = 0.4409056586455051
$ws.Cells.Item(21, 14).Value = 1.785721070759038

# Row 22 (A22 = 20)
$ws.Cells.Item(22, 2).Value = 2.434503014865982
$ws.Cells.Item(22, 3).Value = 0.4762468422443931
$ws.Cells.Item(22, 4).Value = 0.07692236755989512
$ws.Cells.Item(22, 5).Value = 0.05432124290277862
$ws.Cells.Item(22, 6).Value = 2.114726622392283
$ws.Cells.Item(22, 9).Value = 1.440251412531339
$ws.Cells.Item(22, 12).Value = 0.2383207497939281
$ws.Cells.Item(22, 13).Value = 0.4571558660696766
$ws.Cells.Item(22, 14).Value = 1.771735458920524

# Row 23 (A23 = 21)
$ws.Cells.Item(23, 2).Value = 2.377755429431261
$ws.Cells.Item(23, 3).Value = 0.4590406886941878
$ws.Cells.Item(23, 4).Value = 0.07724584522838818
$ws.Cells.Item(23, 5).Value = 0.05405547559004731
$ws.Cells.Item(23, 6).Value = 2.095399434385996
$ws.Cells.Item(23, 9).Value = 1.434940148256771
$ws.Cells.Item(23, 12).Value = 0.2365473672510348
$ws.Cells.Item(23, 13).Value = 0.4484655178082164
$ws.Cells.Item(23, 14).Value = 1.779144968058418

# Row 24 (A24 = 22)
$ws.Cells.Item(24, 2).Value = 2.164663317616487
$ws.Cells.Item(24, 3).Value = 0.39398928115628
$ws.Cells.Item(24, 4).Value = 0.07852238795122801
$ws.Cells.Item(24, 5).Value = 0.0530383595298769
$ws.Cells.Item(24, 6).Value = 2.024516932556253
$ws.Cells.Item(24, 9).Value = 1.416296005725187
$ws.Cells.Item(24, 12).Value = 0.2300374850298965
$ws.Cells.Item(24, 13).Value = 0.4159313989773707
$ws.Cells.Item(24, 14).Value = 1.808410144339867

# Row 25 (A25 = 23)
$ws.Cells.Item(25, 2).Value = 1.938679524494603
$ws.Cells.Item(25, 3).Value = 0.3240391843713155
$ws.Cells.Item(25, 4).Value = 0.08000848297414542
$ws.Cells.Item(25, 5).Value = 0.05191809749380738
$ws.Cells.Item(25, 6).Value = 1.952946163012044
$ws.Cells.Item(25, 9).Value = 1.399299058224926
$ws.Cells.Item(25, 12).Value = 0.2234552268838996
$ws.Cells.Item(25, 13).Value = 0.3816418884016031
$ws.Cells.Item(25, 14).Value = 1.842530103081927
